$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new reference-table rows for "day" (rows 30-36) and "duration" (rows 37-39)
$dayRows = @(
    @("day", 1, "Monday"),
    @("day", 2, "Tuesday"),
    @("day", 3, "Wednesday"),
    @("day", 4, "Thursday"),
    @("day", 5, "Friday"),
    @("day", 6, "Saturday"),
    @("day", 7, "Sunday")
)

$durationRows = @(
    @("duration", 1, "Normal"),
    @("duration", 2, "Separuh Hari"),
    @("duration", 3, "Hujung Minggu")
)

$allRows = $dayRows + $durationRows

$startRow = 30
$r = $startRow
foreach ($row in $allRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Update the view: scroll/zoom and selection to match the edited state
$win = $excel.ActiveWindow
$win.Zoom = 115
$win.ScrollRow = 16
$win.ScrollColumn = 1

$ws.Range("C14").Select()
